$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old sub-header row (Hiver/Eté/Année, (m3/s)/(MW)/(GWh)) - this
# shifts all the plant data rows up by one.
$ws.Rows("2:2").Delete()

# Clear the remaining old header row (E1,G1,I1,J1,K1) formatting+contents so
# we can rebuild it as a single-row header with idx/idx2/Name/Date Start/Date End.
$ws.Range("A1:K1").ClearFormats()
$ws.Range("A1:K1").ClearContents()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

$ws.Range("A2:K2").Select()
